{"js": "// Update the answer cells in the \"two-digit \u00f7 one-digit\" practice table.\n// The table lays its 25 problems out as 5 \"data\" rows (each followed by\n// blank spacer rows) of 5 columns each. Several old answer strings repeat\n// (e.g. \"95\u00f72=47, 1\" and \"69\u00f78=8, 5\" each occur twice), so we must replace\n// the text by (row, column) position rather than by searching for the old\n// value, otherwise a value-based replace could hit the wrong duplicate.\n\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nconst table = tables.items[0];\nconst rows = table.rows;\nrows.load(\"items\");\nawait context.sync();\n\n// Indices (within the 20 <w:tr> rows) of the 5 rows that actually hold\n// problem/answer text; the rows in between are empty spacer rows.\nconst dataRowIndexes = [0, 4, 8, 12, 16];\n\n// New answer text, in row-major order (5 rows x 5 columns), replacing the\n// old text at the same position.\nconst newAnswers = [\n  [\"58\u00f78=7, 2\", \"67\u00f76=11, 1\", \"39\u00f78=4, 7\", \"10\u00f74=2, 2\", \"27\u00f78=3, 3\"],\n  [\"29\u00f75=5, 4\", \"94\u00f73=31, 1\", \"69\u00f78=8, 5\", \"91\u00f76=15, 1\", \"68\u00f76=11, 2\"],\n  [\"70\u00f74=17, 2\", \"83\u00f78=10, 3\", \"57\u00f77=8, 1\", \"49\u00f76=8, 1\", \"27\u00f73=9, 0\"],\n  [\"88\u00f73=29, 1\", \"35\u00f72=17, 1\", \"75\u00f74=18, 3\", \"43\u00f74=10, 3\", \"81\u00f77=11, 4\"],\n  [\"15\u00f74=3, 3\", \"69\u00f72=34, 1\", \"55\u00f79=6, 1\", \"23\u00f72=11, 1\", \"52\u00f79=5, 7\"],\n];\n\n// Load the cells collection for each data row first.\nfor (const rowIndex of dataRowIndexes) {\n  rows.items[rowIndex].cells.load(\"items\");\n}\nawait context.sync();\n\n// Load the first paragraph of each cell body (so we can replace its text\n// in place, preserving the paragraph/run formatting already applied).\nconst rowParagraphs = [];\nfor (const rowIndex of dataRowIndexes) {\n  const cells = rows.items[rowIndex].cells.items;\n  const paragraphs = [];\n  for (const cell of cells) {\n    cell.body.paragraphs.load(\"items\");\n    paragraphs.push(cell.body.paragraphs);\n  }\n  rowParagraphs.push(paragraphs);\n}\nawait context.sync();\n\n// Now replace the text of the first paragraph in each cell with the new\n// answer text, keeping existing run/paragraph formatting.\nfor (let r = 0; r < dataRowIndexes.length; r++) {\n  const paragraphs = rowParagraphs[r];\n  for (let c = 0; c < paragraphs.length; c++) {\n    const firstParagraph = paragraphs[c].items[0];\n    firstParagraph.insertText(newAnswers[r][c], Word.InsertLocation.replace);\n  }\n}\nawait context.sync();\n", "ps1": "# Update the answer cells in the \"two-digit \u00f7 one-digit\" practice table.\n# The table lays its 25 problems out as 5 \"data\" rows (each followed by\n# blank spacer rows) of 5 columns each. Several old answer strings repeat\n# (e.g. \"95\u00f72=47, 1\" and \"69\u00f78=8, 5\" each occur twice), so we replace the\n# text by (row, column) position rather than searching for the old value,\n# which keeps us from accidentally touching the wrong duplicate.\n\n$d = $word.ActiveDocument\n$tbl = $d.Tables.Item(1)\n\n# 1-based row indexes (within the table's 20 rows) that hold problem/answer\n# text; the rows in between are empty spacer rows.\n$dataRows = @(1, 5, 9, 13, 17)\n\n# New answer text, in row-major order (5 rows x 5 columns), replacing the\n# old text at the same position.\n$newAnswers = @(\n    @(\"58\u00f78=7, 2\", \"67\u00f76=11, 1\", \"39\u00f78=4, 7\", \"10\u00f74=2, 2\", \"27\u00f78=3, 3\"),\n    @(\"29\u00f75=5, 4\", \"94\u00f73=31, 1\", \"69\u00f78=8, 5\", \"91\u00f76=15, 1\", \"68\u00f76=11, 2\"),\n    @(\"70\u00f74=17, 2\", \"83\u00f78=10, 3\", \"57\u00f77=8, 1\", \"49\u00f76=8, 1\", \"27\u00f73=9, 0\"),\n    @(\"88\u00f73=29, 1\", \"35\u00f72=17, 1\", \"75\u00f74=18, 3\", \"43\u00f74=10, 3\", \"81\u00f77=11, 4\"),\n    @(\"15\u00f74=3, 3\", \"69\u00f72=34, 1\", \"55\u00f79=6, 1\", \"23\u00f72=11, 1\", \"52\u00f79=5, 7\")\n)\n\nfor ($i = 0; $i -lt $dataRows.Length; $i++) {\n    $r = $dataRows[$i]\n    for ($c = 1; $c -le 5; $c++) {\n        $cell = $tbl.Cell($r, $c)\n        # Range excludes the trailing cell-mark (end-of-cell char) so this\n        # replaces just the paragraph's text, preserving its run/paragraph\n        # formatting (font, size, left alignment, etc.).\n        $rng = $cell.Range\n        $rng.End = $rng.End - 1\n        $rng.Text = $newAnswers[$i][$c - 1]\n    }\n}\n"}
